$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: the "Test R2" (F15) pipeline result is removed (the cell goes away
# entirely rather than just becoming blank) - matches the source XML where
# the <c r="F15"/> element is dropped.
$ws.Range("F15").ClearContents()

# Row 16: newly added pipeline run - "2024-11-20 14:50:24"
$ws.Range("A16").Value = "2024-11-20 14:50:24"
$ws.Range("B16").Value = 0.9961673926144942
$ws.Range("C16").Value = 0.007499415425961405
$ws.Range("D16").Value = 0.0001739766272678393
$ws.Range("E16").Value = 0.01319001998739347
# F16 (Test R2) has no result for this pipeline, same as F15 above - leave blank.
$ws.Range("G16").Value = 0.01775764939020219
$ws.Range("H16").Value = 0.0003776926955150451
$ws.Range("I16").Value = 0.01943431746975039
